# Modification plantType dans l'outil d'importation.
# The "Type vegetale" column (N) used the French label "Légume" to mark
# plant type; it is replaced with the English label "Vegetable".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Rows whose "Type vegetale" (column N) cell currently holds "Légume".
$rows = @(2, 3, 6, 7, 8, 9, 10)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 14)  # column N
    if ($cell.Value() -eq "Légume") {
        $cell.Value = "Vegetable"
    }
}

# Reflect the author's final selection on the sheet (matches the saved
# sheetView selection in the edited workbook).
$ws.Activate()
$appWindow = $excel.ActiveWindow
$appWindow.ScrollColumn = 3
$appWindow.ScrollRow = 4
$ws.Range("N8").Select()

